$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.204.95'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '1.784.21'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'226.54"
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = "'31.88"
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('D10').Value = "'0.0692"
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').Value = '2.041.31'
$ws.Range('D13').Value = "'11.01"
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '1.788.68'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = "'0.624"
$ws.Range('E15').Value = '  +2.44%  '
$ws.Range('D16').Value = '34.181.81'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').Value = "'67.96"
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').Value = '0.0₃0804'
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('D20').Value = "'246.66"
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').Value = "'10.98"
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = "'162.73"
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').Value = "'16.32"
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('E31').Value = '  +2.29%  '
$ws.Range('E32').Value = '  +4.50%  '
$ws.Range('E33').Value = '  +6.94%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '1.445.35'
$ws.Range('E35').Value = '  +4.13%  '
$ws.Range('D36').Value = "'0.655"
$ws.Range('E36').Value = '  +3.05%  '
$ws.Range('D37').Value = "'2.39"
$ws.Range('E37').Value = '  +6.43%  '
$ws.Range('D39').Value = "'1.04"
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').Value = "'80.33"
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('D42').Value = "'0.925"
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('D44').Value = "'13.45"
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').Value = "'6.09"
$ws.Range('E45').Value = '  +3.95%  '
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '0.0₆0135'
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('D49').Value = '1.943.50'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').Value = "'104.58"
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('E51').Value = '  +0.12%  '
